# "change save excel files names"
#
# The report's "report" sheet (projects table) gains a new project column
# (sf_cust) and the "service_tables" sheet renames/reorders the
# san_switch_report_isl -> san_isl_report_tables group (which also flips
# which step-group's export/force_extract flags are turned on), plus two
# isl/trunk rows get their export flags turned off.

$wb = $excel.ActiveWorkbook

# Best-effort: VBA codeName metadata (engine may not persist these, but
# they are harmless to set).
try { $wb.CodeName = "ThisWorkbook" } catch {}

$ws1 = $wb.Worksheets.Item(1)   # "report"
$ws2 = $wb.Worksheets.Item(2)   # "service_tables"

try { $ws1.CodeName = "Sheet1" } catch {}
try { $ws2.CodeName = "Sheet2" } catch {}

# ---------------------------------------------------------------------
# Sheet "service_tables":
#  - rows 16/17 (isl/trunk under san_isl step 11): export_to_excel and
#    force_extract flags flip from 1 to 0.
#  - the "san_fabrics_statistics" (step 17) and "san_switch_report_tables"
#    (step 16) groups swap order: switch_report_tables (rows 33-38) now
#    comes before fabrics_statistics (rows 39-41).
#  - "san_switch_report_isl" (step 18) is renamed to "san_isl_report_tables"
#    and its rows' flags flip from 1 to 0; switch_params_aggregated /
#    fabric_statistics rows flip from 0 to 1.
# ---------------------------------------------------------------------

$ws2.Range("E16").Value = 0
$ws2.Range("F16").Value = 0
$ws2.Range("E17").Value = 0
$ws2.Range("F17").Value = 0

# Row 42: step 18 renamed san_switch_report_isl -> san_isl_report_tables
# (done before the sheet1 edits below so new shared strings land in the
# same order as the authored commit).
$ws2.Range("A42").Value = 18
$ws2.Range("B42").Value = "san_isl_report_tables"
$ws2.Range("D42").Value = "Межкоммутаторные_соединения"
$ws2.Range("E42").Value = 0
$ws2.Range("F42").Value = 0

# ---------------------------------------------------------------------
# Sheet "report": the "value" header moves from C1 to D1, and a new
# project column H is added (sf_cust), mirroring the existing B..G
# project columns.
# ---------------------------------------------------------------------

$ws1.Range("C1").Clear()

$ws1.Range("D1").Font.Bold = $true
$ws1.Range("D1").Interior.Color = $ws1.Range("A1").Interior.Color
$ws1.Range("D1").Value = "value"

$ws1.Range("H4").Value = "C:\Users\vlasenko\Documents\01.CUSTOMERS\Megafon\All SANs\SF\packed_sshow"
$ws1.Range("H3").Value = "C:\Users\vlasenko\Documents\01.CUSTOMERS\Megafon\All SANs\SF"
$ws1.Range("H2").Font.Bold = $ws1.Range("B2").Font.Bold
$ws1.Range("H2").Value = "sf_cust"

# ---------------------------------------------------------------------
# Back to "service_tables" for the remaining row reshuffle.
# ---------------------------------------------------------------------

# Row 33: now step 16 / san_switch_report_tables / Коммутаторы
$ws2.Range("A33").Value = 16
$ws2.Range("B33").Value = "san_switch_report_tables"
$ws2.Range("D33").Value = "Коммутаторы"
$ws2.Range("E33").Value = 1
$ws2.Range("F33").Value = 0

# Row 34: Фабрика
$ws2.Range("D34").Value = "Фабрика"
$ws2.Range("E34").Value = 0
$ws2.Range("F34").Value = 0

# Row 35: Глобальные_параметры_фабрики
$ws2.Range("D35").Value = "Глобальные_параметры_фабрики"
$ws2.Range("E35").Value = 0
$ws2.Range("F35").Value = 0

# Row 36: Параметры_коммутаторов (A36/B36 no longer used - group header
# already written on row 33)
$ws2.Range("A36").ClearContents()
$ws2.Range("B36").ClearContents()
$ws2.Range("D36").Value = "Параметры_коммутаторов"
$ws2.Range("E36").Value = 0
$ws2.Range("F36").Value = 0

# Row 37: Лицензии
$ws2.Range("D37").Value = "Лицензии"
$ws2.Range("E37").Value = 0
$ws2.Range("F37").Value = 0

# Row 38: switch_params_aggregated
$ws2.Range("D38").Value = "switch_params_aggregated"
$ws2.Range("E38").Value = 1
$ws2.Range("F38").Value = 0

# Row 39: now step 17 / san_fabrics_statistics / fabric_statistics
$ws2.Range("A39").Value = 17
$ws2.Range("B39").Value = "san_fabrics_statistics"
$ws2.Range("D39").Value = "fabric_statistics"
$ws2.Range("E39").Value = 1
$ws2.Range("F39").Value = 0

# Row 40: Статистика
$ws2.Range("D40").Value = "Статистика"
$ws2.Range("E40").Value = 0
$ws2.Range("F40").Value = 0

# Row 41: Статистика_Итого
$ws2.Range("D41").Value = "Статистика_Итого"
$ws2.Range("E41").Value = 0
$ws2.Range("F41").Value = 0

# Row 43: Межфабричные_соединения
$ws2.Range("D43").Value = "Межфабричные_соединения"
$ws2.Range("E43").Value = 0
$ws2.Range("F43").Value = 0

# Row 44: isl_aggregated
$ws2.Range("D44").Value = "isl_aggregated"
$ws2.Range("E44").Value = 0
$ws2.Range("F44").Value = 0

# ---------------------------------------------------------------------
# View state: selection / scroll position changes.
# ---------------------------------------------------------------------

$ws1.Activate()
$ws1.Range("A2").Select()

$ws2.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws2.Range("G37:G38").Select()
